$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'275.20"
$ws.Range("E2").Value = "'-1.35%"
$ws.Range("D3").Value = "'27.75"
$ws.Range("E3").Value = "'2.57%"
$ws.Range("D4").Value = "'4.795"
$ws.Range("E4").Value = "'-2.04%"
$ws.Range("D5").Value = "'0.06302"
$ws.Range("E5").Value = "'-0.47%"
$ws.Range("D6").Value = "'6.922"
$ws.Range("E6").Value = "'-0.09%"
$ws.Range("D7").Value = "'1.285"
$ws.Range("E7").Value = "'35.80%"
$ws.Range("D8").Value = "'0.8706"
$ws.Range("E8").Value = "'-1.14%"
$ws.Range("D9").Value = "'0.1522"
$ws.Range("E9").Value = "'3.38%"
$ws.Range("D10").Value = "'0.05010"
$ws.Range("E10").Value = "'-2.37%"
$ws.Range("D11").Value = "'0.07492"
$ws.Range("E11").Value = "'1.78%"
$ws.Range("D12").Value = "'0.02916"
$ws.Range("E12").Value = "'-7.14%"
$ws.Range("D13").Value = "'0.09019"
$ws.Range("E13").Value = "'-0.48%"
$ws.Range("E14").Value = "'0.35%"
$ws.Range("D15").Value = "'0.0006355"
$ws.Range("E15").Value = "'1.43%"
$ws.Range("D16").Value = "'0.006025"
$ws.Range("E16").Value = "'4.53%"
$ws.Range("E17").Value = "'-0.49%"
$ws.Range("D18").Value = "'3.305"
$ws.Range("E18").Value = "'-1.46%"
$ws.Range("D19").Value = "'2.272"
$ws.Range("E19").Value = "'-0.98%"
$ws.Range("D20").Value = "'0.3118"
$ws.Range("E20").Value = "'0.60%"
$ws.Range("D21").Value = "'0.1316"
$ws.Range("E21").Value = "'-1.60%"
$ws.Range("D22").Value = "'3.915"
$ws.Range("E22").Value = "'-0.12%"
$ws.Range("D23").Value = "'0.04383"
$ws.Range("E23").Value = "'1.77%"
$ws.Range("D24").Value = "'0.001169"
$ws.Range("E24").Value = "'-0.69%"
$ws.Range("D25").Value = "'0.003822"
$ws.Range("E25").Value = "'5.79%"
$ws.Range("D26").Value = "'0.0001200"
$ws.Range("E26").Value = "'0.04%"
$ws.Range("D27").Value = "'0.0001651"
$ws.Range("E27").Value = "'-2.20%"
$ws.Range("D40").Value = "'0.04102"
$ws.Range("E40").Value = "'1.56%"
$ws.Range("D41").Value = "'0.007035"
$ws.Range("E41").Value = "'6.42%"
$ws.Range("D42").Value = "'0.1170"
$ws.Range("E42").Value = "'0.54%"
$ws.Range("D43").Value = "'0.002020"
$ws.Range("E43").Value = "'-13.64%"
$ws.Range("D44").Value = "'0.01167"
$ws.Range("E44").Value = "'-6.47%"
$ws.Range("D45").Value = "'0.00005191"
$ws.Range("E45").Value = "'-0.19%"
$ws.Range("D47").Value = "'0.02300"
$ws.Range("E47").Value = "'2.31%"
